# Update recomputed TPM-based NATMI ligand-receptor signalling values
# for Sema6d-Trem2 (rows 2-21, columns E:T) per updated script outputs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 55.91100466666666
$ws.Range("H2").Value = 167.733014
$ws.Range("I2").Value = 0.7311810562391667
$ws.Range("J2").Value = 0.7311810562391668
$ws.Range("M2").Value = 0.4425036666666667
$ws.Range("N2").Value = 1.327511
$ws.Range("O2").Value = 0.00265625830230164
$ws.Range("P2").Value = 0.00265625830230164
$ws.Range("Q2").Value = 24.74082457201711
$ws.Range("R2").Value = 222.667421148154
$ws.Range("S2").Value = 0.001942205751120969
$ws.Range("T2").Value = 0.001942205751120969
# Row 3
$ws.Range("G3").Value = 55.91100466666666
$ws.Range("H3").Value = 167.733014
$ws.Range("I3").Value = 0.7311810562391667
$ws.Range("J3").Value = 0.7311810562391668
$ws.Range("M3").Value = 70.96028133333334
$ws.Range("N3").Value = 212.880844
$ws.Range("O3").Value = 0.4259599425360546
$ws.Range("P3").Value = 0.4259599425360546
$ws.Range("Q3").Value = 3967.46062077598
$ws.Range("R3").Value = 35707.14558698382
$ws.Range("S3").Value = 0.3114538406990872
$ws.Range("T3").Value = 0.3114538406990872
# Row 4
$ws.Range("G4").Value = 55.91100466666666
$ws.Range("H4").Value = 167.733014
$ws.Range("I4").Value = 0.7311810562391667
$ws.Range("J4").Value = 0.7311810562391668
$ws.Range("K4").Value = 1.0
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.080319
$ws.Range("N4").Value = 0.240957
$ws.Range("O4").Value = 0.0004821384016762922
$ws.Range("P4").Value = 0.0004821384016762922
$ws.Range("Q4").Value = 4.490715983822
$ws.Range("R4").Value = 40.416443854398
$ws.Range("S4").Value = 0.0003525304657911349
$ws.Range("T4").Value = 0.000352530465791135
# Row 5
$ws.Range("G5").Value = 55.91100466666666
$ws.Range("H5").Value = 167.733014
$ws.Range("I5").Value = 0.7311810562391667
$ws.Range("J5").Value = 0.7311810562391668
$ws.Range("M5").Value = 95.105991
$ws.Range("N5").Value = 285.317973
$ws.Range("O5").Value = 0.5709016607599675
$ws.Range("P5").Value = 0.5709016607599675
$ws.Range("Q5").Value = 5317.471506628958
$ws.Range("R5").Value = 47857.24355966062
$ws.Range("S5").Value = 0.4174324793231675
$ws.Range("T5").Value = 0.4174324793231676
# Row 6
$ws.Range("I6").Value = 0.1683066942221897
$ws.Range("J6").Value = 0.1683066942221898
$ws.Range("M6").Value = 0.4425036666666667
$ws.Range("N6").Value = 1.327511
$ws.Range("O6").Value = 0.00265625830230164
$ws.Range("P6").Value = 0.00265625830230164
$ws.Range("Q6").Value = 5.694959354479334
$ws.Range("R6").Value = 51.25463419031401
$ws.Range("S6").Value = 0.000447066053860635
$ws.Range("T6").Value = 0.0004470660538606351
# Row 7
$ws.Range("I7").Value = 0.1683066942221897
$ws.Range("J7").Value = 0.1683066942221898
$ws.Range("M7").Value = 70.96028133333334
$ws.Range("N7").Value = 212.880844
$ws.Range("O7").Value = 0.4259599425360546
$ws.Range("P7").Value = 0.4259599425360546
$ws.Range("Q7").Value = 913.2487444000508
$ws.Range("R7").Value = 8219.238699600457
$ws.Range("S7").Value = 0.07169190979931725
$ws.Range("T7").Value = 0.07169190979931726
# Row 8
$ws.Range("I8").Value = 0.1683066942221897
$ws.Range("J8").Value = 0.1683066942221898
$ws.Range("K8").Value = 1.0
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.080319
$ws.Range("N8").Value = 0.240957
$ws.Range("O8").Value = 0.0004821384016762922
$ws.Range("P8").Value = 0.0004821384016762922
$ws.Range("Q8").Value = 1.033694124702
$ws.Range("R8").Value = 9.303247122318
$ws.Range("S8").Value = 0.000081147120543707
$ws.Range("T8").Value = 0.00008114712054370702
# Row 9
$ws.Range("I9").Value = 0.1683066942221897
$ws.Range("J9").Value = 0.1683066942221898
$ws.Range("M9").Value = 95.105991
$ws.Range("N9").Value = 285.317973
$ws.Range("O9").Value = 0.5709016607599675
$ws.Range("P9").Value = 0.5709016607599675
$ws.Range("Q9").Value = 1224.000599119278
$ws.Range("R9").Value = 11016.0053920735
$ws.Range("S9").Value = 0.09608657124846816
$ws.Range("T9").Value = 0.09608657124846817
# Row 10
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 0.6944773333333333
$ws.Range("H10").Value = 2.083432
$ws.Range("I10").Value = 0.00908208809961812
$ws.Range("J10").Value = 0.009082088099618121
$ws.Range("M10").Value = 0.4425036666666667
$ws.Range("N10").Value = 1.327511
$ws.Range("O10").Value = 0.00265625830230164
$ws.Range("P10").Value = 0.00265625830230164
$ws.Range("Q10").Value = 0.3073087664168889
$ws.Range("R10").Value = 2.765778897752
$ws.Range("S10").Value = 0.00002412437191684555
$ws.Range("T10").Value = 0.00002412437191684556
# Row 11
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = 0.6944773333333333
$ws.Range("H11").Value = 2.083432
$ws.Range("I11").Value = 0.00908208809961812
$ws.Range("J11").Value = 0.009082088099618121
$ws.Range("M11").Value = 70.96028133333334
$ws.Range("N11").Value = 212.880844
$ws.Range("O11").Value = 0.4259599425360546
$ws.Range("P11").Value = 0.4259599425360546
$ws.Range("Q11").Value = 49.28030695295644
$ws.Range("R11").Value = 443.522762576608
$ws.Range("S11").Value = 0.003868605725020719
$ws.Range("T11").Value = 0.00386860572502072
# Row 12
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 0.6944773333333333
$ws.Range("H12").Value = 2.083432
$ws.Range("I12").Value = 0.00908208809961812
$ws.Range("J12").Value = 0.009082088099618121
$ws.Range("K12").Value = 1.0
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.080319
$ws.Range("N12").Value = 0.240957
$ws.Range("O12").Value = 0.0004821384016762922
$ws.Range("P12").Value = 0.0004821384016762922
$ws.Range("Q12").Value = 0.055779724936
$ws.Range("R12").Value = 0.5020175244239999
$ws.Range("S12").Value = 0.000004378823440233154
$ws.Range("T12").Value = 0.000004378823440233155
# Row 13
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 0.6944773333333333
$ws.Range("H13").Value = 2.083432
$ws.Range("I13").Value = 0.00908208809961812
$ws.Range("J13").Value = 0.009082088099618121
$ws.Range("M13").Value = 95.105991
$ws.Range("N13").Value = 285.317973
$ws.Range("O13").Value = 0.5709016607599675
$ws.Range("P13").Value = 0.5709016607599675
$ws.Range("Q13").Value = 66.048955013704
$ws.Range("R13").Value = 594.4405951233359
$ws.Range("S13").Value = 0.005184979179240322
$ws.Range("T13").Value = 0.005184979179240322
# Row 14
$ws.Range("G14").Value = 5.954760333333334
$ws.Range("H14").Value = 17.864281
$ws.Range("I14").Value = 0.07787389935372697
$ws.Range("J14").Value = 0.077873899353727
$ws.Range("M14").Value = 0.4425036666666667
$ws.Range("N14").Value = 1.327511
$ws.Range("O14").Value = 0.00265625830230164
$ws.Range("P14").Value = 0.00265625830230164
$ws.Range("Q14").Value = 2.635003281621223
$ws.Range("R14").Value = 23.715029534591
$ws.Range("S14").Value = 0.0002068531916909396
$ws.Range("T14").Value = 0.0002068531916909397
# Row 15
$ws.Range("G15").Value = 5.954760333333334
$ws.Range("H15").Value = 17.864281
$ws.Range("I15").Value = 0.07787389935372697
$ws.Range("J15").Value = 0.077873899353727
$ws.Range("M15").Value = 70.96028133333334
$ws.Range("N15").Value = 212.880844
$ws.Range("O15").Value = 0.4259599425360546
$ws.Range("P15").Value = 0.4259599425360546
$ws.Range("Q15").Value = 422.5514685259072
$ws.Range("R15").Value = 3802.963216733165
$ws.Range("S15").Value = 0.03317116169377204
$ws.Range("T15").Value = 0.03317116169377205
# Row 16
$ws.Range("G16").Value = 5.954760333333334
$ws.Range("H16").Value = 17.864281
$ws.Range("I16").Value = 0.07787389935372697
$ws.Range("J16").Value = 0.077873899353727
$ws.Range("K16").Value = 1.0
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.080319
$ws.Range("N16").Value = 0.240957
$ws.Range("O16").Value = 0.0004821384016762922
$ws.Range("P16").Value = 0.0004821384016762922
$ws.Range("Q16").Value = 0.478280395213
$ws.Range("R16").Value = 4.304523556917
$ws.Range("S16").Value = 0.00003754599736670637
$ws.Range("T16").Value = 0.00003754599736670638
# Row 17
$ws.Range("G17").Value = 5.954760333333334
$ws.Range("H17").Value = 17.864281
$ws.Range("I17").Value = 0.07787389935372697
$ws.Range("J17").Value = 0.077873899353727
$ws.Range("M17").Value = 95.105991
$ws.Range("N17").Value = 285.317973
$ws.Range("O17").Value = 0.5709016607599675
$ws.Range("P17").Value = 0.5709016607599675
$ws.Range("Q17").Value = 566.333382669157
$ws.Range("R17").Value = 5097.000444022414
$ws.Range("S17").Value = 0.04445833847089729
$ws.Range("T17").Value = 0.04445833847089731
# Row 18
$ws.Range("E18").Value = 3.0
$ws.Range("F18").Value = 1.0
$ws.Range("G18").Value = 1.036602666666667
$ws.Range("H18").Value = 3.109808
$ws.Range("I18").Value = 0.01355626208529831
$ws.Range("J18").Value = 0.01355626208529831
$ws.Range("M18").Value = 0.4425036666666667
$ws.Range("N18").Value = 1.327511
$ws.Range("O18").Value = 0.00265625830230164
$ws.Range("P18").Value = 0.00265625830230164
$ws.Range("Q18").Value = 0.4587004808764445
$ws.Range("R18").Value = 4.128304327888
$ws.Range("S18").Value = 0.00003600893371225059
$ws.Range("T18").Value = 0.00003600893371225059
# Row 19
$ws.Range("E19").Value = 3.0
$ws.Range("F19").Value = 1.0
$ws.Range("G19").Value = 1.036602666666667
$ws.Range("H19").Value = 3.109808
$ws.Range("I19").Value = 0.01355626208529831
$ws.Range("J19").Value = 0.01355626208529831
$ws.Range("M19").Value = 70.96028133333334
$ws.Range("N19").Value = 212.880844
$ws.Range("O19").Value = 0.4259599425360546
$ws.Range("P19").Value = 0.4259599425360546
$ws.Range("Q19").Value = 73.55761685755024
$ws.Range("R19").Value = 662.018551717952
$ws.Range("S19").Value = 0.005774424618857364
$ws.Range("T19").Value = 0.005774424618857365
# Row 20
$ws.Range("E20").Value = 3.0
$ws.Range("F20").Value = 1.0
$ws.Range("G20").Value = 1.036602666666667
$ws.Range("H20").Value = 3.109808
$ws.Range("I20").Value = 0.01355626208529831
$ws.Range("J20").Value = 0.01355626208529831
$ws.Range("K20").Value = 1.0
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.080319
$ws.Range("N20").Value = 0.240957
$ws.Range("O20").Value = 0.0004821384016762922
$ws.Range("P20").Value = 0.0004821384016762922
$ws.Range("Q20").Value = 0.08325888958400002
$ws.Range("R20").Value = 0.749330006256
$ws.Range("S20").Value = 0.000006535994534510647
$ws.Range("T20").Value = 0.000006535994534510648
# Row 21
$ws.Range("E21").Value = 3.0
$ws.Range("F21").Value = 1.0
$ws.Range("G21").Value = 1.036602666666667
$ws.Range("H21").Value = 3.109808
$ws.Range("I21").Value = 0.01355626208529831
$ws.Range("J21").Value = 0.01355626208529831
$ws.Range("M21").Value = 95.105991
$ws.Range("N21").Value = 285.317973
$ws.Range("O21").Value = 0.5709016607599675
$ws.Range("P21").Value = 0.5709016607599675
$ws.Range("Q21").Value = 98.58712388657601
$ws.Range("R21").Value = 887.284114979184
$ws.Range("S21").Value = 0.007739292538194186
$ws.Range("T21").Value = 0.007739292538194187
